$wb = $excel.ActiveWorkbook

$wsDelete = $wb.Worksheets.Item("deleteEntity")
$wsRelation = $wb.Worksheets.Item("getRelationById")

# --- deleteEntity sheet (sheet3) ---
# E2: 1234 -> 123456789 (style s="3", no quote-prefix, plain value set is fine)
$wsDelete.Range("E2").Value = 123456789

# E4: 9999 -> 999999, but must keep its quote-prefix style (s="5").
# Setting .Value alone resets the cell to the non-quote-prefix style, so
# after assigning the new value we re-apply the original (quote-prefixed)
# number format/style by pasting formats from a cell that already uses it.
$wsDelete.Range("E4").Value = 999999
$wsDelete.Range("E3").Copy()
$wsDelete.Range("E4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- getRelationById sheet (sheet6) ---
# D5: 9999 -> 999999 (style s="2", no quote-prefix, plain value set is fine)
$wsRelation.Range("D5").Value = 999999

# --- Sheet view / selection changes ---
# getRelationById: selection moves from C5 to D5, and it is no longer the
# active/selected tab.
$wsRelation.Range("D5").Select()

# deleteEntity: selection moves from H8 to F7, and it becomes the
# active/selected tab (activeTab 5 -> 2 in workbook.xml).
$wsDelete.Activate()
$wsDelete.Range("F7").Select()
